$d = $word.ActiveDocument

# wdReplaceAll = 2, wdFindContinue = 1
$replacements = @(
    @("319÷6=53, 1", "716÷2=358, 0"),
    @("233÷2=116, 1", "284÷9=31, 5"),
    @("883÷3=294, 1", "844÷7=120, 4"),
    @("571÷8=71, 3", "642÷5=128, 2"),
    @("163÷7=23, 2", "887÷6=147, 5"),
    @("223÷7=31, 6", "888÷2=444, 0"),
    @("300÷3=100, 0", "409÷3=136, 1"),
    @("511÷7=73, 0", "961÷7=137, 2"),
    @("220÷8=27, 4", "733÷5=146, 3"),
    @("347÷7=49, 4", "588÷2=294, 0"),
    @("855÷7=122, 1", "306÷9=34, 0"),
    @("471÷2=235, 1", "358÷9=39, 7"),
    @("195÷4=48, 3", "181÷8=22, 5"),
    @("804÷6=134, 0", "607÷8=75, 7"),
    @("105÷2=52, 1", "489÷2=244, 1"),
    @("670÷7=95, 5", "787÷6=131, 1"),
    @("737÷3=245, 2", "648÷5=129, 3"),
    @("151÷4=37, 3", "251÷4=62, 3"),
    @("750÷5=150, 0", "547÷2=273, 1"),
    @("417÷6=69, 3", "410÷3=136, 2"),
    @("349÷2=174, 1", "412÷5=82, 2"),
    @("567÷8=70, 7", "398÷7=56, 6"),
    @("400÷9=44, 4", "231÷3=77, 0"),
    @("111÷8=13, 7", "995÷5=199, 0"),
    @("233÷4=58, 1", "755÷4=188, 3"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
